$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(2, 0, 0.0021868999999999999, 45686.686356194259, "TrueTR"),
    @(3, 1, 3.8212373999999998, 45686.686400229286, "TrueTR"),
    @(4, 1, 4.0191001000000002, 45686.68640256011, "SimTR"),
    @(5, 2, 5.0207367999999999, 45686.686414133357, "TrueTR"),
    @(6, 2, 6.0196991000000004, 45686.686425693057, "SimTR"),
    @(7, 3, 7.5409136999999999, 45686.686443273204, "TrueTR"),
    @(8, 3, 8.0368031999999996, 45686.686449009678, "SimTR"),
    @(9, 4, 10.0511208, 45686.68647232365, "SimTR"),
    @(10, 4, 10.056268899999999, 45686.686472382593, "TrueTR"),
    @(11, 5, 12.005111899999999, 45686.686494941023, "SimTR"),
    @(12, 5, 12.5400432, 45686.686501129057, "TrueTR"),
    @(13, 6, 14.0394325, 45686.686518486342, "SimTR"),
    @(14, 6, 15.021493700000001, 45686.686529850544, "TrueTR"),
    @(15, 7, 15.072672600000001, 45686.686530442712, "TrueTR"),
    @(16, 7, 16.0204369, 45686.686541413197, "SimTR"),
    @(17, 8, 17.5398979, 45686.686558997913, "TrueTR"),
    @(18, 8, 18.020530999999998, 45686.686564561111, "SimTR"),
    @(19, 9, 20.021532400000002, 45686.686587720877, "SimTR"),
    @(20, 9, 20.023060699999998, 45686.686587737437, "TrueTR"),
    @(21, 10, 20.0711206, 45686.686588294811, "TrueTR"),
    @(22, 10, 22.039170500000001, 45686.68661107299, "SimTR"),
    @(23, 11, 22.539384500000001, 45686.686616864034, "TrueTR"),
    @(24, 12, 22.555964299999999, 45686.686617054584, "TrueTR"),
    @(25, 11, 24.020783099999999, 45686.686634009166, "SimTR"),
    @(26, 13, 25.0259578, 45686.686645642294, "TrueTR"),
    @(27, 12, 26.002044300000001, 45686.686656942671, "SimTR"),
    @(28, 14, 27.515249600000001, 45686.686674453755, "TrueTR"),
    @(29, 15, 27.549562399999999, 45686.686674851604, "TrueTR"),
    @(30, 13, 28.00001, 45686.686680063627, "SimTR"),
    @(31, 14, 30.000041499999998, 45686.686703212843, "SimTR"),
    @(32, 16, 30.042865899999999, 45686.686703707986, "TrueTR"),
    @(33, 17, 30.058609300000001, 45686.686703891108, "TrueTR"),
    @(34, 15, 32.000104999999998, 45686.686726361862, "SimTR"),
    @(35, 18, 32.538333899999998, 45686.686732590962, "TrueTR"),
    @(36, 19, 32.552061600000002, 45686.686732750059, "TrueTR"),
    @(37, 16, 34.000054400000003, 45686.686749510925, "SimTR"),
    @(38, 17, 36.4229828, 45686.686777552648, "SimTR"),
    @(39, 20, 36.424269199999998, 45686.686777565548, "TrueTR"),
    @(40, 21, 37.536952599999999, 45686.686790445543, "TrueTR"),
    @(41, 22, 37.570059499999999, 45686.686790829408, "TrueTR"),
    @(42, 18, 38.0031319, 45686.686795841153, "SimTR"),
    @(43, 19, 40.003137299999999, 45686.6868189897, "SimTR"),
    @(44, 23, 40.037225399999997, 45686.686819385905, "TrueTR"),
    @(45, 24, 40.070220499999998, 45686.686819765535, "TrueTR"),
    @(46, 20, 42.003271300000002, 45686.686842138864, "SimTR"),
    @(47, 25, 42.536633700000003, 45686.68684831178, "TrueTR"),
    @(48, 26, 42.570053799999997, 45686.686848698773, "TrueTR"),
    @(49, 21, 44.004962300000003, 45686.686865307573, "SimTR"),
    @(50, 27, 45.0381134, 45686.686877263914, "TrueTR"),
    @(51, 28, 45.072978999999997, 45686.686877668719, "TrueTR"),
    @(52, 22, 46.004399599999999, 45686.686888449185, "SimTR"),
    @(53, 29, 47.538185900000002, 45686.686906199931, "TrueTR"),
    @(54, 30, 47.571886999999997, 45686.686906590243, "TrueTR"),
    @(55, 23, 48.003859499999997, 45686.686911590587, "SimTR"),
    @(56, 24, 50.003642800000002, 45686.686934736295, "SimTR"),
    @(57, 31, 50.070648200000001, 45686.686935511236, "TrueTR"),
    @(58, 25, 52.003482599999998, 45686.686957881626, "SimTR")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}